$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("M2").Value = 0.05057900000000001
$ws.Range("N2").Value = 0.151737
$ws.Range("O2").Value = 0.01400296657613869
$ws.Range("P2").Value = 0.01400296657613869
$ws.Range("Q2").Value = 4.197433917299001
$ws.Range("R2").Value = 37.77690525569101
$ws.Range("S2").Value = 0.006286637607177472
$ws.Range("T2").Value = 0.006286637607177472
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("O3").Value = 0.146324388539341
$ws.Range("P3").Value = 0.146324388539341
$ws.Range("Q3").Value = 43.86120241332667
$ws.Range("R3").Value = 394.75082171994
$ws.Range("S3").Value = 0.06569239445348501
$ws.Range("T3").Value = 0.06569239445348501
$ws.Range("G4").Value = 82.98768099999999
$ws.Range("H4").Value = 248.963043
$ws.Range("I4").Value = 0.4489504115427952
$ws.Range("J4").Value = 0.4489504115427952
$ws.Range("O4").Value = 0.8396726448845202
$ws.Range("P4").Value = 0.8396726448845202
$ws.Range("Q4").Value = 251.6945548575547
$ws.Range("R4").Value = 2265.250993717992
$ws.Range("S4").Value = 0.3769713794821327
$ws.Range("T4").Value = 0.3769713794821327
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("M5").Value = 0.05057900000000001
$ws.Range("N5").Value = 0.151737
$ws.Range("O5").Value = 0.01400296657613869
$ws.Range("P5").Value = 0.01400296657613869
$ws.Range("Q5").Value = 3.193587867890667
$ws.Range("R5").Value = 28.742290811016
$ws.Range("S5").Value = 0.004783143698668746
$ws.Range("T5").Value = 0.004783143698668746
$ws.Range("G6").Value = 63.14058933333333
$ws.Range("I6").Value = 0.3415807409566563
$ws.Range("J6").Value = 0.3415807409566563
$ws.Range("O6").Value = 0.146324388539341
$ws.Range("P6").Value = 0.146324388539341
$ws.Range("R6").Value = 300.34336690544
$ws.Range("S6").Value = 0.04998159305729776
$ws.Range("T6").Value = 0.04998159305729776
$ws.Range("G7").Value = 63.14058933333333
$ws.Range("I7").Value = 0.3415807409566563
$ws.Range("J7").Value = 0.3415807409566563
$ws.Range("O7").Value = 0.8396726448845202
$ws.Range("P7").Value = 0.8396726448845202
$ws.Range("S7").Value = 0.2868160042006898
$ws.Range("T7").Value = 0.2868160042006898
$ws.Range("I8").Value = 0.2094688475005485
$ws.Range("J8").Value = 0.2094688475005485
$ws.Range("M8").Value = 0.05057900000000001
$ws.Range("N8").Value = 0.151737
$ws.Range("O8").Value = 0.01400296657613869
$ws.Range("P8").Value = 0.01400296657613869
$ws.Range("Q8").Value = 1.958415946418
$ws.Range("R8").Value = 17.625743517762
$ws.Range("S8").Value = 0.002933185270292473
$ws.Range("T8").Value = 0.002933185270292473
$ws.Range("I9").Value = 0.2094688475005485
$ws.Range("J9").Value = 0.2094688475005485
$ws.Range("O9").Value = 0.146324388539341
$ws.Range("P9").Value = 0.146324388539341
$ws.Range("S9").Value = 0.03065040102855822
$ws.Range("T9").Value = 0.03065040102855822
$ws.Range("I10").Value = 0.2094688475005485
$ws.Range("J10").Value = 0.2094688475005485
$ws.Range("O10").Value = 0.8396726448845202
$ws.Range("P10").Value = 0.8396726448845202
$ws.Range("S10").Value = 0.1758852612016978
$ws.Range("T10").Value = 0.1758852612016978